# Add an optional "Priority" column (before the existing "Status" column, i.e.
# the old column AA) to both the "Screen Print Designs" and "Embroidery Designs"
# sheets, and populate it for a handful of rows that already had a Status set.

$wb = $excel.ActiveWorkbook

# --- Screen Print Designs (sheet 1) ---------------------------------------
$ws1 = $wb.Worksheets.Item("Screen Print Designs")
$ws1.Activate()

# Insert a new column at AA (27); everything from the old AA onward shifts
# right by one (AA -> AB, AB -> AC, ... AJ -> AK) and all formulas/validations
# that reference those columns are adjusted automatically.
$ws1.Columns.Item(27).Insert()

$ws1.Cells.Item(1, 27).Value = "Priority"
$ws1.Cells.Item(2, 27).Value = 3
$ws1.Cells.Item(6, 27).Value = 1
$ws1.Cells.Item(10, 27).Value = 2

$null = $ws1.Range("AA12").Select()

# --- Embroidery Designs (sheet 2) -----------------------------------------
$ws2 = $wb.Worksheets.Item("Embroidery Designs")
$ws2.Activate()

$ws2.Columns.Item(27).Insert()

$ws2.Cells.Item(1, 27).Value = "Priority"
$ws2.Cells.Item(3, 27).Value = 1
$ws2.Cells.Item(7, 27).Value = 2
$ws2.Cells.Item(6, 27).Value = 3

$null = $ws2.Range("AA9").Select()

$ws1.Activate()
